$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RESOUCE AVILABILITY")

# Update the "Status" column (C) for rows 2-6, matching the look of the
# existing YES/NO cells in column D (fill color copied from there).
$ws.Range("C2:C5").Value = "YES"
$ws.Range("C2:C5").Interior.Color = $ws.Range("D3").Interior.Color

$ws.Range("C6").Value = "NO"
$ws.Range("C6").Interior.Color = $ws.Range("D2").Interior.Color

# Move the active selection in the bottom-right frozen pane to G4.
$ws.Range("G4").Select()
